# Add new weapons list
# - Items sheet: 5 new rows (Garrote, Poisoned Wind Globe, Scythe, Spike Fist, Things Catcher)
# - Feats sheet: 3 new rows (Poisoned Wind Fumigator, Spike Fist Brutalist, Things Wrangler)
# - Re-point the active tab from MagicItems to Feats

$wb = $excel.ActiveWorkbook

$newWeaponsUrl = "https://editor.gmbinder.com/documents/edit/-N8HgfCb1-XMt8weFOvO"

# ---------------------------------------------------------------------------
# Items sheet (sheet5.xml) -- new rows 6-10
# ---------------------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("Items")

$itemsData = @(
    @("Garrote", "Weapon"),
    @("Poisoned Wind Globe", "Item"),
    @("Scythe", "Weapon"),
    @("Spike Fist", "Weapon"),
    @("Things Catcher", "Weapon")
)

$r = 6
foreach ($row in $itemsData) {
    $wsItems.Range("A$r").Value = $row[0]
    $wsItems.Range("B$r").Value = $row[1]
    $wsItems.Range("C$r").Value = "New Weapons"
    $wsItems.Range("D$r").Value = "Playtest Ready"
    $wsItems.Range("E$r").Value = "Not Released"

    # Hyperlink the item-name (column C) cell, reusing the same target doc as
    # the existing "New Weapons" hyperlinks in this sheet.
    $wsItems.Hyperlinks.Add($wsItems.Range("C$r"), $newWeaponsUrl)

    # Hyperlinks.Add() mints a brand-new cell style; restore the shared
    # "Hyperlink"-style (s="3") already used by C2:C5 via a style-only paste.
    $wsItems.Range("C2").Copy()
    $wsItems.Range("C$r").PasteSpecial(-4122)

    $r++
}

# Conditional formatting on column D needs to keep covering the new rows.
$wsItems.Range("D1:D10").FormatConditions.Item(1).ModifyAppliesToRange($wsItems.Range("D1:D10"))

# Column A got wider to fit the new (longer) item names.
$wsItems.Columns.Item(1).EntireColumn.AutoFit()

$wsItems.Range("D11").Select()

# ---------------------------------------------------------------------------
# Feats sheet (sheet7.xml) -- new rows 18-20
# ---------------------------------------------------------------------------
$wsFeats = $wb.Worksheets.Item("Feats")

$featsData = @(
    "Poisoned Wind Fumigator",
    "Spike Fist Brutalist",
    "Things Wrangler"
)

$r = 18
foreach ($name in $featsData) {
    $wsFeats.Range("A$r").Value = $name
    $wsFeats.Range("B$r").Value = "Weapon"
    $wsFeats.Range("C$r").Value = "None"
    $wsFeats.Range("D$r").Value = "No"
    $wsFeats.Range("E$r").Value = "New Weapons"
    $wsFeats.Range("F$r").Value = "Playtest Ready"
    $wsFeats.Range("G$r").Value = "Not Released"
    $r++
}

# E18 gets its own single-cell hyperlink; E19:E20 share one merged hyperlink
# range (mirrors how the existing multi-row "Independent Feats" links work).
$wsFeats.Hyperlinks.Add($wsFeats.Range("E18"), $newWeaponsUrl)
$wsFeats.Hyperlinks.Add($wsFeats.Range("E19:E20"), $newWeaponsUrl, "", "", "New Weapons")

# Restore the shared hyperlink style (s="3") on the three new E cells.
$wsFeats.Range("E17").Copy()
$wsFeats.Range("E18").PasteSpecial(-4122)
$wsFeats.Range("E19:E20").PasteSpecial(-4122)

# Column A got wider to fit the new (longer) feat names.
$wsFeats.Columns.Item(1).EntireColumn.AutoFit()

$wsFeats.Range("C23").Select()

# Feats becomes the active/selected tab (was MagicItems before).
$wsFeats.Activate()
